# TestCase_A14 / TestCase_A15: only the "A Suite" (row 2) should keep
# running; B/C/D/E/F Suite (rows 3-7) switch their Runmode from Y to N,
# and the visible selection moves down to reflect the now-disabled rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3:C7").Value = "N"

$ws.Range("C5:C7").Select()
